$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# A7 value changes from "Yo" to "Gabriel" (shared string replaced)
$ws.Range("A7").Value = "Gabriel"

# B7 and C7 get the "highlighted" style used elsewhere in the sheet (same as C6),
# i.e. green font color, keep the thin border already present.
$ws.Range("B7:C7").Font.Color = $ws.Range("C6").Font.Color

# Update the selected cell / active cell shown when the workbook is opened
$ws.Range("C12").Select()

$wb.Save()
